$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 21 de Septiembre de 2020 a las 19:04"

# Row 4: Estados Unidos
$ws.Range("B4").Value = 7016028
$ws.Range("C4").Value = 11260
$ws.Range("D4").Value = 4269757
$ws.Range("E4").Value = 2542056
$ws.Range("G4").Value = 97
$ws.Range("H4").Value = 204215

# Row 6: Brasil
$ws.Range("B6").Value = 4547150
$ws.Range("C6").Value = 2521
$ws.Range("E6").Value = 558926
$ws.Range("G6").Value = 102
$ws.Range("H6").Value = 136997

# Row 22: Turquia
$ws.Range("B22").Value = 304610
$ws.Range("C22").Value = 1743
$ws.Range("D22").Value = 268435
$ws.Range("E22").Value = 28601
$ws.Range("G22").Value = 68
$ws.Range("H22").Value = 7574

# Row 31: Ecuador
$ws.Range("B31").Value = 126711
$ws.Range("C31").Value = 292
$ws.Range("E31").Value = 12764
$ws.Range("G31").Value = 5
$ws.Range("H31").Value = 11095

# Row 51: Etiopia
$ws.Range("A51").Value = "Etiopia"
$ws.Range("B51").Value = 69709
$ws.Range("C51").Value = 889
$ws.Range("D51").Value = 28634
$ws.Range("E51").Value = 39967
$ws.Range("G51").Value = 12
$ws.Range("H51").Value = 1108

# Row 52: Portugal
$ws.Range("A52").Value = "Portugal"
$ws.Range("B52").Value = 69200
$ws.Range("C52").Value = 623
$ws.Range("D52").Value = 45736
$ws.Range("E52").Value = 21544
$ws.Range("G52").Value = 8
$ws.Range("H52").Value = 1920

# Row 57: Singapur
$ws.Range("B57").Value = 57606
$ws.Range("C57").Value = 30
$ws.Range("D57").Value = 57241
$ws.Range("E57").Value = 338

# Row 61: Chequia
$ws.Range("A61").Value = "Chequia"
$ws.Range("B61").Value = 50071
$ws.Range("C61").Value = 781
$ws.Range("D61").Value = 25336
$ws.Range("E61").Value = 24214
$ws.Range("G61").Value = 18
$ws.Range("H61").Value = 521

# Row 62: Argelia
$ws.Range("A62").Value = "Argelia"
$ws.Range("B62").Value = 50023
$ws.Range("C62").Value = 197
$ws.Range("D62").Value = 35180
$ws.Range("E62").Value = 13164
$ws.Range("G62").Value = 7
$ws.Range("H62").Value = 1679

# Row 75: Libano
$ws.Range("B75").Value = 29987
$ws.Range("C75").Value = 684
$ws.Range("D75").Value = 12507
$ws.Range("E75").Value = 17173
$ws.Range("G75").Value = 10
$ws.Range("H75").Value = 307

# Row 88: Grecia
$ws.Range("B88").Value = 15595
$ws.Range("C88").Value = 453
$ws.Range("E88").Value = 5262
$ws.Range("G88").Value = 6
$ws.Range("H88").Value = 344

# Row 110: Mozambique
$ws.Range("B110").Value = 6912
$ws.Range("C110").Value = 141
$ws.Range("D110").Value = 3738
$ws.Range("E110").Value = 3130
$ws.Range("G110").Value = 1
$ws.Range("H110").Value = 44

# Row 140: Reunion
$ws.Range("A140").Value = "Reunion"
$ws.Range("B140").Value = 3415
$ws.Range("C140").Value = 73
$ws.Range("D140").Value = 2482
$ws.Range("E140").Value = 918
$ws.Range("H140").Value = 15

# Row 141: Bahamas
$ws.Range("A141").Value = "Bahamas"
$ws.Range("B141").Value = 3370
$ws.Range("C141").Value = 55
$ws.Range("D141").Value = 1689
$ws.Range("E141").Value = 1607
$ws.Range("H141").Value = 74

# Row 142: Sri Lanka
$ws.Range("A142").Value = "Sri Lanka"
$ws.Range("B142").Value = 3298
$ws.Range("C142").Value = 11
$ws.Range("D142").Value = 3100
$ws.Range("E142").Value = 185
$ws.Range("H142").Value = 13

# Row 145: Malta
$ws.Range("E145").Value = 674
$ws.Range("G145").Value = 3
$ws.Range("H145").Value = 23

# Row 146: Sudan del Sur
$ws.Range("B146").Value = 2649
$ws.Range("C146").Value = 7
$ws.Range("E146").Value = 1310

# Row 157: Principado de Andorra
$ws.Range("A157").Value = "Principado de Andorra"
$ws.Range("B157").Value = 1681
$ws.Range("C157").Value = 117
$ws.Range("D157").Value = 1199
$ws.Range("E157").Value = 429
$ws.Range("H157").Value = 53

# Row 158: Togo
$ws.Range("A158").Value = "Togo"
$ws.Range("B158").Value = 1666
$ws.Range("C158").Value = 0
$ws.Range("D158").Value = 1269
$ws.Range("E158").Value = 356
$ws.Range("G158").Value = 0
$ws.Range("H158").Value = 41

# Row 159: Belice
$ws.Range("A159").Value = "Belice"
$ws.Range("B159").Value = 1627
$ws.Range("C159").Value = 21
$ws.Range("D159").Value = 918
$ws.Range("E159").Value = 688
$ws.Range("G159").Value = 1
$ws.Range("H159").Value = 21

# Row 160: Republica de Chipre
$ws.Range("A160").Value = "Republica de Chipre"
$ws.Range("B160").Value = 1603
$ws.Range("C160").Value = 3
$ws.Range("D160").Value = 1369
$ws.Range("E160").Value = 212
$ws.Range("H160").Value = 22

# Row 177: Burundi
$ws.Range("B177").Value = 474
$ws.Range("C177").Value = 1
$ws.Range("E177").Value = 11

# Row 186: Curazao
$ws.Range("A186").Value = "Curazao"
$ws.Range("B186").Value = 282
$ws.Range("C186").Value = 14
$ws.Range("D186").Value = 102

# Row 187: Camboya
$ws.Range("A187").Value = "Camboya"
$ws.Range("B187").Value = 275
$ws.Range("D187").Value = 274
$ws.Range("E187").Value = 1
$ws.Range("H187").Value = 0

# Row 195: Liechtenstein
$ws.Range("B195").Value = 114
$ws.Range("C195").Value = 1
$ws.Range("D195").Value = 110

# Row 204: Santa Lucia
$ws.Range("A204").Value = "Santa Lucia"

# Row 205: Timor Oriental
$ws.Range("A205").Value = "Timor Oriental"

# Row 214: Montserrat
$ws.Range("A214").Value = "Montserrat"
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1

# Row 215: Islas Malvinas
$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("D215").Value = 13
$ws.Range("H215").Value = 0
